# Adds 10 new rating rows (rows 7-16) to the Rating sheet, each holding a
# single elapsed-time string in the column matching its difficulty
# (A = Easy, B = Normal, C = Hard), mirroring OpenRandomCell writing a new
# finish time into the ratings table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimes = @(
    @{ Row = 7;  Col = "A"; Value = "00:17:60" },
    @{ Row = 8;  Col = "B"; Value = "00:44:94" },
    @{ Row = 9;  Col = "A"; Value = "00:07:54" },
    @{ Row = 10; Col = "A"; Value = "03:32:84" },
    @{ Row = 11; Col = "C"; Value = "01:14:28" },
    @{ Row = 12; Col = "B"; Value = "01:30:01" },
    @{ Row = 13; Col = "A"; Value = "00:05:39" },
    @{ Row = 14; Col = "A"; Value = "00:06:16" },
    @{ Row = 15; Col = "B"; Value = "00:34:93" },
    @{ Row = 16; Col = "C"; Value = "01:06:20" }
)

foreach ($entry in $newTimes) {
    $addr = "$($entry.Col)$($entry.Row)"
    $ws.Range($addr).Value = $entry.Value
}
